$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) ---
$volCell = $ws.Range("A8")
$volCell.Characters(21, 2).Text = "36"
$volCell.Characters(21, 2).Font.Name = "Andale WT"
$volCell.Characters(21, 2).Font.Size = 10

$dateCell = $ws.Range("C9")
$dateCell.Characters(47, 8).Text = "9/10/2023"
$dateCell.Characters(47, 9).Font.Name = "Andale WT"
$dateCell.Characters(47, 9).Font.Size = 10
$dateCell.Characters(27, 9).Text = "9/4/2023"
$dateCell.Characters(27, 8).Font.Name = "Andale WT"
$dateCell.Characters(27, 8).Font.Size = 10

# --- Data table updates (rows 15-30) ---
$ws.Range("C14").Copy($ws.Range("F15"))
$ws.Range("H15").Value = -100
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -42.857142857142
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -28.571428571428
$ws.Range("I16").Value = 141
$ws.Range("J16").Value = 173
$ws.Range("K16").Value = -18.497109826589
$ws.Range("L16").Value = 25.892857142857
$ws.Range("M16").Value = 50
$ws.Range("N16").Value = -84.984025559105
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 50
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 30
$ws.Range("I17").Value = 131
$ws.Range("J17").Value = 122
$ws.Range("K17").Value = 7.377049180327
$ws.Range("L17").Value = 33.673469387755
$ws.Range("M17").Value = 101.538461538462
$ws.Range("N17").Value = -34.825870646766
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -40
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -15
$ws.Range("I18").Value = 174
$ws.Range("J18").Value = 171
$ws.Range("K18").Value = 1.754385964912
$ws.Range("L18").Value = 16.778523489932
$ws.Range("M18").Value = -1.694915254237
$ws.Range("N18").Value = -90.965732087227
$ws.Range("C19").Value = 28
$ws.Range("D19").Value = 35
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 141
$ws.Range("G19").Value = 157
$ws.Range("H19").Value = -10.191082802547
$ws.Range("I19").Value = 1184
$ws.Range("J19").Value = 1198
$ws.Range("K19").Value = -1.168614357262
$ws.Range("L19").Value = 60
$ws.Range("M19").Value = 37.037037037037
$ws.Range("N19").Value = -53.164556962025
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -71.428571428571
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -28.571428571428
$ws.Range("I20").Value = 116
$ws.Range("J20").Value = 132
$ws.Range("K20").Value = -12.121212121212
$ws.Range("L20").Value = 2.654867256637
$ws.Range("M20").Value = 110.909090909091
$ws.Range("N20").Value = -95.040615647712
$ws.Range("C21").Value = 40
$ws.Range("E21").Value = -28.571428571428
$ws.Range("F21").Value = 197
$ws.Range("G21").Value = 223
$ws.Range("H21").Value = -11.659192825112
$ws.Range("I21").Value = 1758
$ws.Range("J21").Value = 1807
$ws.Range("K21").Value = -2.711676812396
$ws.Range("L21").Value = 43.862520458265
$ws.Range("M21").Value = 39.082278481012
$ws.Range("N21").Value = -77.903469079939
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("C14").Copy($ws.Range("F22"))
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -100
$ws.Range("M22").Value = 60
$ws.Range("C14").Copy($ws.Range("C23"))
$ws.Range("I14").Copy($ws.Range("D23"))
$ws.Range("D23").Value = 3
$ws.Range("K14").Copy($ws.Range("E23"))
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -66.666666666666
$ws.Range("I23").Value = 16
$ws.Range("J23").Value = 23
$ws.Range("K23").Value = -30.434782608695
$ws.Range("L23").Value = -40.740740740740
$ws.Range("M23").Value = -11.111111111111
$ws.Range("C24").Value = 63
$ws.Range("D24").Value = 81
$ws.Range("E24").Value = -22.222222222222
$ws.Range("F24").Value = 238
$ws.Range("G24").Value = 374
$ws.Range("H24").Value = -36.363636363636
$ws.Range("I24").Value = 2267
$ws.Range("J24").Value = 2745
$ws.Range("K24").Value = -17.413479052823
$ws.Range("L24").Value = 33.904311872415
$ws.Range("M24").Value = 95.937770095073
$ws.Range("C25").Value = 5
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 20
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = -23.076923076923
$ws.Range("I25").Value = 221
$ws.Range("J25").Value = 255
$ws.Range("K25").Value = -13.333333333333
$ws.Range("L25").Value = -2.212389380530
$ws.Range("M25").Value = -9.795918367346
$ws.Range("C14").Copy($ws.Range("F26"))
$ws.Range("H26").Value = -100
$ws.Range("I14").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = -22.222222222222
$ws.Range("I27").Value = 69
$ws.Range("J27").Value = 59
$ws.Range("K27").Value = 16.949152542372
$ws.Range("L27").Value = 7.8125
$ws.Range("C14").Copy($ws.Range("D30"))
$ws.Range("E14").Copy($ws.Range("E30"))
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 100
